# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# D = Price, E = Volume(1h) text columns; both are stored as plain text in the
# sheet (e.g. "26.041.76" is a display string, not a number). For D-column
# values that look like a single plain decimal (e.g. "216.55"), Excel's COM
# layer would otherwise auto-coerce the assigned string into a real number,
# so those cells are briefly forced to Text format, written, and then
# restored to the default "Normal" style so no stray formatting is left
# behind (matching the source file, where these cells carry no style index).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.041.76"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.646.40"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  +0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.508"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "1.874.90"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "1.648.22"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.546"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").Value = "26.057.32"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("E25").Value = "  +5.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("E35").Value = "  +2.31%  "
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").Value = "1.129.92"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.541"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.798"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "1.784.16"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("E45").Value = "  +4.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.40%  "
